$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the "scenario1" rows (4 & 5) and "scenario2" rows (9 & 10) date
# values from the 2022 test dates to the new 2023 test dates. These cells are
# stored as quote-prefixed text (so they display literally instead of being
# reformatted as serial dates) - prefix the literal with an apostrophe so the
# cell keeps being treated as text, matching the original authoring style.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value2 = "'01/01/2023"
$ws.Range("C5").Value2 = "'01/05/2023"
$ws.Range("C9").Value2 = "'01/01/2023"
$ws.Range("C10").Value2 = "'01/05/2023"
$ws.Range("D9").Value2 = "'01/03/2023"
$ws.Range("D10").Value2 = "'01/13/2023"

# ---------------------------------------------------------------------------
# Rows 2/3/7/8 column C (and D on rows 7/8) used a left-over "applied but
# empty" alignment style; clear that stray formatting so the cells fall back
# to the default (unstyled) look, matching the cleaned-up workbook.
# ---------------------------------------------------------------------------
$ws.Range("C2").ClearFormats()
$ws.Range("C3").ClearFormats()
$ws.Range("C7").ClearFormats()
$ws.Range("C8").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()

# D2 and D3 were empty placeholder cells carrying that same stray style with
# no content - remove them outright so no cell record remains.
$ws.Range("D2").Clear()
$ws.Range("D3").Clear()

# ---------------------------------------------------------------------------
# Move the active selection to C9 (was D9).
# ---------------------------------------------------------------------------
$ws.Range("C9").Select() | Out-Null
